$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (row, new Price text, new Volume(1h) text) for every changed data row.
# Price/Volume columns hold plain text (not real numbers) in this sheet,
# so every Price write below goes through a temporary Text number-format
# to stop Excel auto-coercing numeric-looking strings (e.g. "235.23") into
# real numbers -- keeping the cell text-typed like the original.
$changes = @(
    ,@(2, "37.491.97", "  +2.58%  ")
    ,@(3, "2.077.11", "  +3.67%  ")
    ,@(4, $null, "  +0.06%  ")
    ,@(5, "235.23", "  -0.29%  ")
    ,@(6, "0.617", "  +3.10%  ")
    ,@(7, "58.47", "  +6.20%  ")
    ,@(8, $null, "  +0.04%  ")
    ,@(9, $null, "  +4.00%  ")
    ,@(10, "59.14", "  +1.58%  ")
    ,@(11, "0.0764", "  +2.31%  ")
    ,@(12, $null, "  +4.07%  ")
    ,@(13, "2.383.80", "  +3.82%  ")
    ,@(14, "14.64", "  +2.83%  ")
    ,@(15, "21.07", "  +3.21%  ")
    ,@(16, "0.780", "  +2.94%  ")
    ,@(17, "5.20", "  +2.25%  ")
    ,@(18, "2.092.66", "  +4.51%  ")
    ,@(19, "37.700.74", "  +3.32%  ")
    ,@(20, "6.24", "  +17.96%  ")
    ,@(21, "70.37", "  +3.69%  ")
    ,@(22, $null, "  +1.49%  ")
    ,@(23, "227.06", "  +2.40%  ")
    ,@(24, $null, "  -0.13%  ")
    ,@(26, $null, "  +0.62%  ")
    ,@(27, "166.50", "  +2.05%  ")
    ,@(28, $null, "  +10.03%  ")
    ,@(29, "9.01", "  +3.88%  ")
    ,@(30, "19.31", "  +2.71%  ")
    ,@(31, "0.127", "  -0.22%  ")
    ,@(32, $null, "  +2.16%  ")
    ,@(33, "4.53", "  +3.13%  ")
    ,@(34, $null, "  +2.94%  ")
    ,@(35, $null, "  +7.12%  ")
    ,@(36, "4.56", "  +6.87%  ")
    ,@(37, $null, "  +0.08%  ")
    ,@(38, "3.36", "  +0.26%  ")
    ,@(39, "1.77", "  +0.62%  ")
    ,@(40, "5.88", "  +3.56%  ")
    ,@(41, "4.66", "  +22.05%  ")
    ,@(42, $null, "  -1.15%  ")
    ,@(43, "0.0953", "  +1.93%  ")
    ,@(44, $null, "  +7.63%  ")
    ,@(45, "96.13", "  +7.30%  ")
    ,@(46, "1.456.42", "  +0.16%  ")
    ,@(47, $null, "  +4.94%  ")
    ,@(48, "15.88", "  +4.30%  ")
    ,@(49, $null, "  +4.06%  ")
    ,@(50, "7.28", "  +6.02%  ")
    ,@(51, $null, "  +1.88%  ")
)

foreach ($change in $changes) {
    $r = $change[0]
    $dVal = $change[1]
    $eVal = $change[2]
    if ($dVal -ne $null) {
        $dCell = $ws.Range("D$r")
        $dCell.NumberFormat = "@"
        $dCell.Value = $dVal
        $dCell.ClearFormats()
    }
    if ($eVal -ne $null) {
        $ws.Range("E$r").Value = $eVal
    }
}
